$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Paragraph 3 ("Introducción:") - drop the _GoBack bookmark at its end.
# ---------------------------------------------------------------------------
$p3xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="000D7B83" w:rsidRDefault="00207081" w:rsidP="000D7B83"><w:pPr><w:spacing w:after="120" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Introducción:</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(3).Range.InsertXML($p3xml)

# ---------------------------------------------------------------------------
# Paragraph 4 ("La aplicación está destinada...") - merge the first four runs
# (keeping the leading tab intact) into a single run, leaving the trailing
# "para el aprendizaje..." run untouched/separate.
# ---------------------------------------------------------------------------
$p4xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="000D7B83" w:rsidRDefault="000D7B83" w:rsidP="00986851"><w:pPr><w:spacing w:after="120" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:tab/><w:t xml:space="preserve">La aplicación está destinada a alumnos y profesores de Primer Ciclo de Educación Primaria. Se trata de un Prototipo vertical de alta fidelidad </w:t></w:r><w:r><w:t>para el aprendizaje de Conocimiento del Medio.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(4).Range.InsertXML($p4xml)

# ---------------------------------------------------------------------------
# Paragraph 14 ("Rol de Toma de medidas") - append the Mattia Rosselli blurb.
# ---------------------------------------------------------------------------
$p14xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00F242A7" w:rsidRDefault="00F242A7" w:rsidP="005A7B56"><w:pPr><w:spacing w:after="120" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Rol de Toma de medidas</w:t></w:r><w:r><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Mattia</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Rosselli</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. Su papel será observar al usuario de la aplicación, ver el tiempo que tarda en realizar las tareas y manejar el programa de captura de pantalla, así como el de conteo de clics.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(14).Range.InsertXML($p14xml)

# ---------------------------------------------------------------------------
# Paragraph 15 (was "Rol de Coordinador") becomes "Rol de Relación con el
# usuario: Jorge Justo Vergés. ..."
# ---------------------------------------------------------------------------
$p15xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00F242A7" w:rsidRDefault="00F242A7" w:rsidP="005A7B56"><w:pPr><w:spacing w:after="120" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t>Rol de Relación con el usuario</w:t></w:r><w:r><w:t xml:space="preserve">: Jorge Justo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Vergés</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. Hará las entrevistas a los usuarios de la aplicación y responderá las dudas de los mismos salvo en las fases donde no se deba para no influir en el usuario.</w:t></w:r></w:p>
'@
$d.Paragraphs.Item(15).Range.InsertXML($p15xml)

# ---------------------------------------------------------------------------
# Paragraph 16 (was "Rol de Relación con el usuario") becomes "Rol de
# Coordinador: Héctor Fernández Matellanes. ..." and gains the _GoBack
# bookmark that used to sit after "Introducción:".
# ---------------------------------------------------------------------------
$p16xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00010075" w:rsidRDefault="00393D60" w:rsidP="005A7B56"><w:pPr><w:spacing w:after="120" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Rol de Coordinador: Héctor Fernández </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Matellanes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. Se encargará de cambiar de fase durante la ejecución de la evaluación con el usuario, controlará los tiempos y reforzará cualquiera de los roles anteriores dependiendo de las necesidades.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$d.Paragraphs.Item(16).Range.InsertXML($p16xml)
